$d = $word.ActiveDocument

function Find-ParaIndex($doc, $prefix) {
    $cnt = $doc.Paragraphs.Count
    for ($i = 1; $i -le $cnt; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        if ($t.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# Step A: swap the order of the two paragraphs in the "Quy trình Mở
# chuyến" section -- "Sau khi thu thập ... tài xế, xe ..." used to
# come before "Nhân viên bắt đầu tiến hành lập lịch ..."; now it
# comes after.
# ------------------------------------------------------------------
$idxSauKhi = Find-ParaIndex $d "Sau khi thu thập thông tin đầy đủ về tài xế"
$pSauKhi = $d.Paragraphs.Item($idxSauKhi)
$rngSauKhi = $d.Range($pSauKhi.Range.Start, $pSauKhi.Range.End)
$rngSauKhi.Cut()

$pAfterCut = $d.Paragraphs.Item($idxSauKhi)
$pasteRange = $d.Range($pAfterCut.Range.End, $pAfterCut.Range.End)
$pasteRange.Paste()

# ------------------------------------------------------------------
# Step B: rework the "Sau khi ..." paragraph's wording & split it
# into three runs, with the middle run ("tiến hành phân công xong")
# tagged as English (US) -- matching an existing en-US run elsewhere
# in the document -- and the _GoBack bookmark sitting right after it.
# ------------------------------------------------------------------
$idxSauKhi = Find-ParaIndex $d "Sau khi thu thập thông tin đầy đủ về tài xế"
$pSauKhi = $d.Paragraphs.Item($idxSauKhi)
$paraStart = $pSauKhi.Range.Start

$run1 = "Sau khi "
$run2 = "tiến hành phân công xong"
$run3 = ", thông tin được ghi nhận và trình bày với ban giám đốc quyết định cuối cùng trước khi triển khai. "

$fullRange = $d.Range($paraStart, $pSauKhi.Range.End - 1)
$fullRange.Text = $run1 + $run2 + $run3

Write-Host "Rewrote paragraph text:" $d.Paragraphs.Item($idxSauKhi).Range.Text

# Borrow the "en-US" character formatting from an existing run of the
# same length elsewhere in the document ("Một chuyến xe được thay "),
# then overwrite just the text of that now-formatted run so the
# language tag sticks while the words change.
$idxDonorPara = Find-ParaIndex $d "Một chuyến xe được thay đổi"
$pDonor = $d.Paragraphs.Item($idxDonorPara)
$donorRange = $d.Range($pDonor.Range.Start, $pDonor.Range.Start + $run2.Length)
Write-Host "Donor text:" $donorRange.Text
$donorFormatted = $donorRange.FormattedText

$run2Start = $paraStart + $run1.Length
$run2End = $run2Start + $run2.Length
$run2Range = $d.Range($run2Start, $run2End)
$run2Range.FormattedText = $donorFormatted

$run2Range2 = $d.Range($run2Start, $run2End)
$run2Range2.Text = $run2

Write-Host "After run2 restyle:" $d.Paragraphs.Item($idxSauKhi).Range.Text

# ------------------------------------------------------------------
# Step C: move the _GoBack bookmark -- delete it from the "Chuyến mới
# được triển khai để xe chạy." paragraph in "Quy trình Sửa chuyến",
# and add it (zero-length) right between run 2 and run 3 above.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$bookmarkPos = $run2End
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

Write-Host "Done."
for ($i = 100; $i -le 106; $i++) {
    Write-Host $i ":" $d.Paragraphs.Item($i).Range.Text
}
